$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2800
$ws.Range("I40").Value = 2250
$ws.Range("J40").Value = 2937.5
$ws.Range("K40").Value = 2250
$ws.Range("L40").Value = 2937.5
$ws.Range("M40").Value = -2075
$ws.Range("N40").Value = -3287.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3181.575
$ws.Range("I64").Value = 3006.52
$ws.Range("J64").Value = 3473.3333
$ws.Range("K64").Value = 3006.52
$ws.Range("L64").Value = 3473.3333
$ws.Range("M64").Value = -2758.52
$ws.Range("N64").Value = -3969.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3181.575
$ws.Range("I67").Value = 3006.52
$ws.Range("J67").Value = 3473.3333
$ws.Range("K67").Value = 3006.52
$ws.Range("L67").Value = 3473.3333
$ws.Range("M67").Value = -2148.52
$ws.Range("N67").Value = -5189.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1839.6296
$ws.Range("I137").Value = 1366.6842
$ws.Range("J137").Value = 2962.875
$ws.Range("K137").Value = 4100.0526
$ws.Range("L137").Value = 8888.625
$ws.Range("M137").Value = -1550.0526
$ws.Range("N137").Value = -13988.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 57384.383
$ws.Range("J140").Value = 57384.383
$ws.Range("L140").Value = 57384.383
$ws.Range("N140").Value = -67744.383

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2335.0417
$ws.Range("I61").Value = 2082.9375
$ws.Range("J61").Value = 2839.25
$ws.Range("K61").Value = 2082.9375
$ws.Range("L61").Value = 2839.25
$ws.Range("M61").Value = -1870.9375
$ws.Range("N61").Value = -3263.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2373.75
$ws.Range("J63").Value = 2247.5
$ws.Range("L63").Value = 2247.5
$ws.Range("N63").Value = -3619.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2373.75
$ws.Range("J66").Value = 2247.5
$ws.Range("L66").Value = 11237.5
$ws.Range("N66").Value = -18101.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 935.5294
$ws.Range("I74").Value = 971.93335
$ws.Range("J74").Value = 662.5
$ws.Range("K74").Value = 971.93335
$ws.Range("L74").Value = 662.5
$ws.Range("M74").Value = -97.93335000000002
$ws.Range("N74").Value = -2410.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 935.5294
$ws.Range("I77").Value = 971.93335
$ws.Range("J77").Value = 662.5
$ws.Range("K77").Value = 4859.66675
$ws.Range("L77").Value = 3312.5
$ws.Range("M77").Value = -491.6667500000003
$ws.Range("N77").Value = -12048.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 15016.667
$ws.Range("J92").Value = 15016.667
$ws.Range("L92").Value = 15016.667
$ws.Range("N92").Value = -20008.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 150041470
$ws.Range("J125").Value = 150041470
$ws.Range("L125").Value = 150041470
$ws.Range("N125").Value = -150051310

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 219248.66
$ws.Range("I132").Value = 278856.72
$ws.Range("K132").Value = 836570.1599999999
$ws.Range("M132").Value = -834040.1599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2335.0417
$ws.Range("I136").Value = 2082.9375
$ws.Range("J136").Value = 2839.25
$ws.Range("K136").Value = 6248.8125
$ws.Range("L136").Value = 8517.75
$ws.Range("M136").Value = -3698.8125
$ws.Range("N136").Value = -13617.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3308.3845
$ws.Range("I105").Value = 2456.875
$ws.Range("J105").Value = 4670.8
$ws.Range("K105").Value = 2456.875
$ws.Range("L105").Value = 4670.8
$ws.Range("M105").Value = -709.875
$ws.Range("N105").Value = -8164.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 73503.734
$ws.Range("I134").Value = 92892.63
$ws.Range("K134").Value = 278677.89
$ws.Range("M134").Value = -276142.89

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1809.174
$ws.Range("I31").Value = 1396.4166
$ws.Range("J31").Value = 2259.4546
$ws.Range("K31").Value = 1396.4166
$ws.Range("L31").Value = 2259.4546
$ws.Range("M31").Value = -1101.4166
$ws.Range("N31").Value = -2849.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1809.174
$ws.Range("I34").Value = 1396.4166
$ws.Range("J34").Value = 2259.4546
$ws.Range("K34").Value = 1396.4166
$ws.Range("L34").Value = 2259.4546
$ws.Range("M34").Value = -1194.4166
$ws.Range("N34").Value = -2663.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1456.4242
$ws.Range("I58").Value = 1340
$ws.Range("J58").Value = 1766.8889
$ws.Range("K58").Value = 1340
$ws.Range("L58").Value = 1766.8889
$ws.Range("M58").Value = -1137
$ws.Range("N58").Value = -2172.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3622.2354
$ws.Range("I132").Value = 2969.6667
$ws.Range("J132").Value = 5188.4
$ws.Range("K132").Value = 8909.000100000001
$ws.Range("L132").Value = 15565.2
$ws.Range("M132").Value = -6379.000100000001
$ws.Range("N132").Value = -20625.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4990.147
$ws.Range("I134").Value = 6586.9546
$ws.Range("K134").Value = 19760.8638
$ws.Range("M134").Value = -17225.8638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1456.4242
$ws.Range("I136").Value = 1340
$ws.Range("J136").Value = 1766.8889
$ws.Range("K136").Value = 4020
$ws.Range("L136").Value = 5300.6667
$ws.Range("M136").Value = -1470
$ws.Range("N136").Value = -10400.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 732.6667
$ws.Range("I98").Value = 799
$ws.Range("J98").Value = 600
$ws.Range("K98").Value = 2397
$ws.Range("L98").Value = 1800
$ws.Range("M98").Value = -899
$ws.Range("N98").Value = -4796

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 586.5357
$ws.Range("I113").Value = 572.1111
$ws.Range("J113").Value = 593.3684
$ws.Range("K113").Value = 1716.3333
$ws.Range("L113").Value = 1780.1052
$ws.Range("M113").Value = 453.6667000000002
$ws.Range("N113").Value = -6120.1052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 44245
$ws.Range("I122").Value = 209.6875
$ws.Range("J122").Value = 52943.332
$ws.Range("K122").Value = 1887.1875
$ws.Range("L122").Value = 476489.988
$ws.Range("M122").Value = 562.8125
$ws.Range("N122").Value = -481389.988

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5965
$ws.Range("I133").Value = 4450
$ws.Range("J133").Value = 9500
$ws.Range("K133").Value = 13350
$ws.Range("L133").Value = 28500
$ws.Range("M133").Value = -8290
$ws.Range("N133").Value = -38620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2326.1052
$ws.Range("I132").Value = 1957
$ws.Range("J132").Value = 3359.6
$ws.Range("K132").Value = 5871
$ws.Range("L132").Value = 10078.8
$ws.Range("M132").Value = -3341
$ws.Range("N132").Value = -15138.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2617.6365
$ws.Range("I132").Value = 2288
$ws.Range("J132").Value = 2892.3333
$ws.Range("K132").Value = 6864
$ws.Range("L132").Value = 8676.999899999999
$ws.Range("M132").Value = -4334
$ws.Range("N132").Value = -13736.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2143.9285
$ws.Range("I136").Value = 1601.5
$ws.Range("K136").Value = 4804.5
$ws.Range("M136").Value = -2254.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2749.35
$ws.Range("I132").Value = 2362.5667
$ws.Range("J132").Value = 3909.7
$ws.Range("K132").Value = 7087.7001
$ws.Range("L132").Value = 11729.1
$ws.Range("M132").Value = -4557.7001
$ws.Range("N132").Value = -16789.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1522.8462
$ws.Range("I136").Value = 1273.0605
$ws.Range("J136").Value = 2896.6667
$ws.Range("K136").Value = 3819.1815
$ws.Range("L136").Value = 8690.000100000001
$ws.Range("M136").Value = -1269.1815
$ws.Range("N136").Value = -13790.0001
